# Adding DESI as an acceptable enumeration for ims-v2 ms_source element
$wb = $excel.ActiveWorkbook

# 1. Insert "DESI" into the "ms_source list" sheet, right before "nanoDESI" (row 7),
#    pushing "nanoDESI" down to row 8.
$msSourceWs = $wb.Worksheets.Item("ms_source list")
$msSourceWs.Rows.Item(7).Insert()
$msSourceWs.Cells.Item(7, 1).Value = "DESI"

# 2. Update the data validation on the main sheet (column Q, ms_source) so its
#    list range grows from $A$1:$A$7 to $A$1:$A$8 to include the new entry.
$mainWs = $wb.Worksheets.Item(1)
$validation = $mainWs.Range("Q2:Q1048576").Validation
$validation.Modify(3, 1, 1, "'ms_source list'!`$A`$1:`$A`$8")
